$wb = $excel.ActiveWorkbook

# 1. Insert a new worksheet "data_collection_mode list" right after "ion_mobility list"
#    (pushes lc_length_unit/.../resolution_y_unit list sheets down by one position)
$afterSheet = $wb.Worksheets.Item("ion_mobility list")
$newSheet = $wb.Worksheets.Add($null, $afterSheet)
$newSheet.Name = "data_collection_mode list"
$newSheet.Range("A1").Value = "DDA"
$newSheet.Range("A2").Value = "DIA"
$newSheet.Range("A3").Value = "MRM"
$newSheet.Range("A4").Value = "PRM"

# 2. Update the comment on Y1 (data_collection_mode column) in "Export as TSV"
#    to mention the new MRM / PRM modes
$ws = $wb.Worksheets.Item("Export as TSV")
$cell = $ws.Range("Y1")
[void]$cell.Comment.Text("Mode of data collection in tandem MS assays. Either DDA (Data-dependent acquisition), DIA (Data-independent acquisition), MRM (multiple reaction monitoring), or PRM (parallel reaction monitoring).")

# 3. Add data validation for column Y referencing the new data_collection_mode list
$range = $ws.Range("Y2:Y1048576")
$range.Validation.Add(3, 1, 1, "'data_collection_mode list'!`$A`$1:`$A`$4")
$range.Validation.ErrorTitle = "Value must come from list"
$range.Validation.ErrorMessage = "Value must be one of: DDA / DIA / MRM / PRM."
$range.Validation.ShowInput = $true
$range.Validation.ShowError = $true
